$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Delete()
$tr.Text = "Anomaly detection in graphs - past, present and future."
